$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filling")

# Check-in / Check-out datetime values used throughout the sheet
$checkIn = 42380.333333333336
$checkOut = 42380.708333333336

# Fill in Check_In / Check_Out for existing employee rows 4-10 (B and C columns)
for ($r = 4; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = $checkIn
    $ws.Cells.Item($r, 3).Value = $checkOut
}

# Add employee PAL057 in row 11
$ws.Range("A11").Value = "PAL057 (Bika Alif)"
$ws.Cells.Item(11, 2).Value = $checkIn
$ws.Cells.Item(11, 3).Value = $checkOut

# Add employee PAL058 in row 12
$ws.Range("A12").Value = "PAL058 (Dewa Rahman)"
$ws.Cells.Item(12, 2).Value = $checkIn
$ws.Cells.Item(12, 3).Value = $checkOut

# Extend the employee dropdown (data validation) on column A down to row 12
$ws.Range("A9:A12").Validation.Delete()
$ws.Range("A9:A12").Validation.Add(3, 1, 1, "allEmployee")

# Move the selection to reflect where the user ended up editing
$ws.Range("C4:C12").Select() | Out-Null
